$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 64
$ws.Range("H64").Value = 4499.4
$ws.Range("I64").Value = 4133.2666
$ws.Range("J64").Value = 5048.6
$ws.Range("K64").Value = 4133.2666
$ws.Range("L64").Value = 5048.6
$ws.Range("M64").Value = -3885.2666
$ws.Range("N64").Value = -5544.6

# Row 67
$ws.Range("H67").Value = 4499.4
$ws.Range("I67").Value = 4133.2666
$ws.Range("J67").Value = 5048.6
$ws.Range("K67").Value = 4133.2666
$ws.Range("L67").Value = 5048.6
$ws.Range("M67").Value = -3275.2666
$ws.Range("N67").Value = -6764.6

# Row 98
$ws.Range("H98").Value = 606.9474
$ws.Range("I98").Value = 579.06665
$ws.Range("J98").Value = 711.5
$ws.Range("K98").Value = 579.06665
$ws.Range("L98").Value = 711.5
$ws.Range("M98").Value = 918.93335
$ws.Range("N98").Value = -3707.5

# Row 100
$ws.Range("H100").Value = 4907.0312
$ws.Range("I100").Value = 2123.5
$ws.Range("J100").Value = 6172.273
$ws.Range("K100").Value = 2123.5
$ws.Range("L100").Value = 6172.273
$ws.Range("M100").Value = -1582.5
$ws.Range("N100").Value = -7254.273

# Row 106
$ws.Range("H106").Value = 3855.7896
$ws.Range("I106").Value = 2069.0908
$ws.Range("K106").Value = 2069.0908
$ws.Range("M106").Value = -1438.0908

# Row 113
$ws.Range("H113").Value = 3268.6365
$ws.Range("I113").Value = 2577.2727
$ws.Range("J113").Value = 3960
$ws.Range("K113").Value = 2577.2727
$ws.Range("L113").Value = 3960
$ws.Range("M113").Value = 676.7273
$ws.Range("N113").Value = -10468

# Row 116
$ws.Range("H116").Value = 65568.06
$ws.Range("I116").Value = 78691.21000000001
$ws.Range("J116").Value = 4326.6665
$ws.Range("K116").Value = 78691.21000000001
$ws.Range("L116").Value = 4326.6665
$ws.Range("M116").Value = -75249.21000000001
$ws.Range("N116").Value = -11210.6665

# Row 122
$ws.Range("H122").Value = 606.9474
$ws.Range("I122").Value = 579.06665
$ws.Range("J122").Value = 711.5
$ws.Range("K122").Value = 1737.19995
$ws.Range("L122").Value = 2134.5
$ws.Range("M122").Value = 712.8000500000001
$ws.Range("N122").Value = -7034.5

# Row 132
$ws.Range("H132").Value = 3346.6538
$ws.Range("I132").Value = 2150.8125
$ws.Range("K132").Value = 6452.4375
$ws.Range("M132").Value = -3922.4375

# Row 138
$ws.Range("H138").Value = 1923.44
$ws.Range("I138").Value = 745.9722
$ws.Range("J138").Value = 3010.3333
$ws.Range("K138").Value = 2237.9166
$ws.Range("L138").Value = 9030.999899999999
$ws.Range("M138").Value = 2902.0834
$ws.Range("N138").Value = -19310.9999

$ws = $wb.Worksheets.Item("ARM")
# Row 61
$ws.Range("H61").Value = 437617.6
$ws.Range("I61").Value = 359527
$ws.Range("J61").Value = 559091.9
$ws.Range("K61").Value = 359527
$ws.Range("L61").Value = 559091.9
$ws.Range("M61").Value = -359315
$ws.Range("N61").Value = -559515.9

# Row 92
$ws.Range("H92").Value = 33000
$ws.Range("J92").Value = 33000
$ws.Range("L92").Value = 33000
$ws.Range("N92").Value = -37992

# Row 102
$ws.Range("H102").Value = 13120
$ws.Range("I102").Value = 1866.6666
$ws.Range("J102").Value = 30000
$ws.Range("K102").Value = 1866.6666
$ws.Range("L102").Value = 30000
$ws.Range("M102").Value = -244.6666
$ws.Range("N102").Value = -33244

# Row 136
$ws.Range("H136").Value = 437617.6
$ws.Range("I136").Value = 359527
$ws.Range("J136").Value = 559091.9
$ws.Range("K136").Value = 1078581
$ws.Range("L136").Value = 1677275.7
$ws.Range("M136").Value = -1076031
$ws.Range("N136").Value = -1682375.7

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 2453.9348
$ws.Range("I31").Value = 1465.6562
$ws.Range("J31").Value = 4712.857
$ws.Range("K31").Value = 1465.6562
$ws.Range("L31").Value = 4712.857
$ws.Range("M31").Value = -1170.6562
$ws.Range("N31").Value = -5302.857

# Row 34
$ws.Range("H34").Value = 2453.9348
$ws.Range("I34").Value = 1465.6562
$ws.Range("J34").Value = 4712.857
$ws.Range("K34").Value = 1465.6562
$ws.Range("L34").Value = 4712.857
$ws.Range("M34").Value = -1263.6562
$ws.Range("N34").Value = -5116.857

# Row 105
$ws.Range("H105").Value = 940.42554
$ws.Range("I105").Value = 883.5897
$ws.Range("J105").Value = 1217.5
$ws.Range("K105").Value = 883.5897
$ws.Range("L105").Value = 1217.5
$ws.Range("M105").Value = 863.4103
$ws.Range("N105").Value = -4711.5

# Row 134
$ws.Range("H134").Value = 2245.1875
$ws.Range("I134").Value = 1190.1111
$ws.Range("J134").Value = 3601.7144
$ws.Range("K134").Value = 3570.3333
$ws.Range("L134").Value = 10805.1432
$ws.Range("M134").Value = -1035.3333
$ws.Range("N134").Value = -15875.1432

$ws = $wb.Worksheets.Item("CUL")
# Row 45
$ws.Range("H45").Value = 1000
$ws.Range("I45").Value = 0
$ws.Range("J45").Value = 1000
$ws.Range("K45").Value = 0
$ws.Range("L45").Value = 3000
$ws.Range("M45").ClearContents()
$ws.Range("N45").Value = -4064

# Row 124
$ws.Range("H124").Value = 2876.25
$ws.Range("I124").Value = 0
$ws.Range("J124").Value = 2876.25
$ws.Range("K124").Value = 0
$ws.Range("L124").Value = 8628.75
$ws.Range("M124").ClearContents()
$ws.Range("N124").Value = -18448.75

# Row 130
$ws.Range("H130").Value = 127866.625
$ws.Range("I130").Value = 0
$ws.Range("J130").Value = 127866.625
$ws.Range("K130").Value = 0
$ws.Range("L130").Value = 383599.875
$ws.Range("M130").ClearContents()
$ws.Range("N130").Value = -393639.875

# Row 131
$ws.Range("H131").Value = 2111.612
$ws.Range("I131").Value = 2847.5
$ws.Range("J131").Value = 1951.0546
$ws.Range("K131").Value = 8542.5
$ws.Range("L131").Value = 5853.1638
$ws.Range("M131").Value = -3502.5
$ws.Range("N131").Value = -15933.1638

# Row 133
$ws.Range("H133").Value = 3069
$ws.Range("I133").Value = 1398.75
$ws.Range("J133").Value = 9750
$ws.Range("K133").Value = 4196.25
$ws.Range("L133").Value = 29250
$ws.Range("M133").Value = 863.75
$ws.Range("N133").Value = -39370

$ws = $wb.Worksheets.Item("LTW")
# Row 136
$ws.Range("H136").Value = 5204.857
$ws.Range("I136").Value = 3109.2222
$ws.Range("J136").Value = 7423.7646
$ws.Range("K136").Value = 9327.6666
$ws.Range("L136").Value = 22271.2938
$ws.Range("M136").Value = -6777.6666
$ws.Range("N136").Value = -27371.2938

$ws = $wb.Worksheets.Item("WVR")
# Row 132
$ws.Range("H132").Value = 23257632
$ws.Range("I132").Value = 33334470
$ws.Range("J132").Value = 3395.2307
$ws.Range("K132").Value = 100003410
$ws.Range("L132").Value = 10185.6921
$ws.Range("M132").Value = -100000880
$ws.Range("N132").Value = -15245.6921
